$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p007r_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p007r_1</id>", 2)
$d.Content.Find.Execute("<id>p007r_a2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p007r_2</id>", 2)
$d.Content.Find.Execute("<id>p007r_a3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p007r_3</id>", 2)
$d.Content.Find.Execute("<id>p007r_a4</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p007r_4</id>", 2)
